# Add a date/time value (with a custom number format) into a new column F
# on the worksheet, next to the existing A:D data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell F1 holds a date-time serial value, formatted as dd/mm/yyyy hh:mm AM/PM
$ws.Range("F1").Value = 36892.521
$ws.Range("F1").NumberFormat = "dd/mm/yyyy hh:mm AM/PM"
